$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Nppa", "Npr2", "ECs", 1, 0.3333333333333333, 0.110293, 0.330879, 0.1397687944778318, 0.1397687944778318, 3, 1, 19.84200033333333, 59.52600099999999, 0.4658919963347263, 0.4658919963347263, 2.188433742764333, 19.695903684879, 0.06511716268457511, 0.06511716268457513)
    ,@("ECs", "Nppa", "Npr2", "FAPs", 1, 0.3333333333333333, 0.110293, 0.330879, 0.1397687944778318, 0.1397687944778318, 3, 1, 12.91523933333333, 38.745718, 0.3032510097300563, 0.3032510097300563, 1.424460491791333, 12.820144426122, 0.0423850280541552, 0.04238502805415521)
    ,@("ECs", "Nppa", "Npr2", "MuSCs", 1, 0.3333333333333333, 0.110293, 0.330879, 0.1397687944778318, 0.1397687944778318, 3, 1, 9.805585333333333, 29.416756, 0.2302360472448257, 0.2302360472448257, 1.081487423169333, 9.733386808523999, 0.0321798147687504, 0.0321798147687504)
    ,@("ECs", "Nppa", "Npr2", "Resolving-Mac", 1, 0.3333333333333333, 0.110293, 0.330879, 0.1397687944778318, 0.1397687944778318, 1, 0.3333333333333333, 0.02644566666666667, 0.079337, 0.0006209466903917867, 0.0006209466903917867, 0.002916771913666666, 0.026250947223, [double]"8.678897035105948E-05", [double]"8.678897035105948E-05")
    ,@("FAPs", "Nppa", "Npr2", "ECs", 1, 0.3333333333333333, 0.03729466666666666, 0.111884, 0.04726166302895539, 0.04726166302895539, 3, 1, 19.84200033333333, 59.52600099999999, 0.4658919963347263, 0.4658919963347263, 0.7400007884315554, 6.660007095883999, 0.02201883053865915, 0.02201883053865916)
    ,@("FAPs", "Nppa", "Npr2", "FAPs", 1, 0.3333333333333333, 0.03729466666666666, 0.111884, 0.04726166302895539, 0.04726166302895539, 3, 1, 12.91523933333333, 38.745718, 0.3032510097300563, 0.3032510097300563, 0.4816695458568888, 4.335025912711999, 0.01433214703505239, 0.01433214703505239)
    ,@("FAPs", "Nppa", "Npr2", "MuSCs", 1, 0.3333333333333333, 0.03729466666666666, 0.111884, 0.04726166302895539, 0.04726166302895539, 3, 1, 9.805585333333333, 29.416756, 0.2302360472448257, 0.2302360472448257, 0.3656960364782222, 3.291264328304, 0.0108813384820036, 0.0108813384820036)
    ,@("FAPs", "Nppa", "Npr2", "Resolving-Mac", 1, 0.3333333333333333, 0.03729466666666666, 0.111884, 0.04726166302895539, 0.04726166302895539, 1, 0.3333333333333333, 0.02644566666666667, 0.079337, 0.0006209466903917867, 0.0006209466903917867, 0.0009862823231111112, 0.008876540908, [double]"2.934697324024172E-05", [double]"2.934697324024172E-05")
    ,@("MuSCs", "Nppa", "Npr2", "ECs", 2, 0.6666666666666666, 0.3940343333333334, 1.182103, 0.4993399740044802, 0.4993399740044802, 3, 1, 19.84200033333333, 59.52600099999999, 0.4658919963347263, 0.4658919963347263, 7.818429373344777, 70.365864360103, 0.2326384973386776, 0.2326384973386776)
    ,@("MuSCs", "Nppa", "Npr2", "FAPs", 2, 0.6666666666666666, 0.3940343333333334, 1.182103, 0.4993399740044802, 0.4993399740044802, 3, 1, 12.91523933333333, 38.745718, 0.3032510097300563, 0.3032510097300563, 5.089047720550445, 45.801429484954, 0.1514253513154387, 0.1514253513154387)
    ,@("MuSCs", "Nppa", "Npr2", "MuSCs", 2, 0.6666666666666666, 0.3940343333333334, 1.182103, 0.4993399740044802, 0.4993399740044802, 3, 1, 9.805585333333333, 29.416756, 0.2302360472448257, 0.2302360472448257, 3.863737279763111, 34.773635517868, 0.1149660618461255, 0.1149660618461255)
    ,@("MuSCs", "Nppa", "Npr2", "Resolving-Mac", 2, 0.6666666666666666, 0.3940343333333334, 1.182103, 0.4993399740044802, 0.4993399740044802, 1, 0.3333333333333333, 0.02644566666666667, 0.079337, 0.0006209466903917867, 0.0006209466903917867, 0.01042050063455556, 0.09378450571100001, 0.0003100635042384028, 0.0003100635042384028)
    ,@("Resolving-Mac", "Nppa", "Npr2", "ECs", 1, 0.3333333333333333, 0.2474883333333333, 0.742465, 0.3136295684887327, 0.3136295684887327, 3, 1, 19.84200033333333, 59.52600099999999, 0.4658919963347263, 0.4658919963347263, 4.91066359249611, 44.195972332465, 0.1461175057728144, 0.1461175057728145)
    ,@("Resolving-Mac", "Nppa", "Npr2", "FAPs", 1, 0.3333333333333333, 0.2474883333333333, 0.742465, 0.3136295684887327, 0.3136295684887327, 3, 1, 12.91523933333333, 38.745718, 0.3032510097300563, 0.3032510097300563, 3.196371057207778, 28.76733951487, 0.09510848332541003, 0.09510848332541004)
    ,@("Resolving-Mac", "Nppa", "Npr2", "MuSCs", 1, 0.3333333333333333, 0.2474883333333333, 0.742465, 0.3136295684887327, 0.3136295684887327, 3, 1, 9.805585333333333, 29.416756, 0.2302360472448257, 0.2302360472448257, 2.426767971504444, 21.84091174354, 0.07220883214794614, 0.07220883214794614)
    ,@("Resolving-Mac", "Nppa", "Npr2", "Resolving-Mac", 1, 0.3333333333333333, 0.2474883333333333, 0.742465, 0.3136295684887327, 0.3136295684887327, 1, 0.3333333333333333, 0.02644566666666667, 0.079337, 0.0006209466903917867, 0.0006209466903917867, 0.006544993967222223, 0.05890494570500001, 0.0001947472425620828, 0.0001947472425620828)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
